# Generate Report for Archive
# - Update the localization status from "Ready for handoff" to "In Translation"
#   everywhere it appears (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4 all share the
#   same shared-string entry).
# - Narrow the affected "Status" columns to match the shorter new text
#   (mirrors the autofit-driven column-width shrink seen in the source diff).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# --- Overview sheet: per-locale status columns E (zh-cn) and F (de-de) ---
$wsOverview.Range("E2:F4").Value = $newStatus

# --- zh-cn / de-de sheets: Status column C ---
$wsZhCn.Range("C2:C4").Value = $newStatus
$wsDeDe.Range("C2:C4").Value = $newStatus

# --- Shrink the Status columns to fit the new, shorter text ---
$wsOverview.Columns.Item(5).ColumnWidth = 12.67
$wsOverview.Columns.Item(6).ColumnWidth = 12.67
$wsZhCn.Columns.Item(3).ColumnWidth = 12.67
$wsDeDe.Columns.Item(3).ColumnWidth = 12.67
